$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.282089710235596
$ws.Range("B1").Value = 4.74585485458374
$ws.Range("C1").Value = 3.917789459228516
$ws.Range("D1").Value = 4.710554599761963
$ws.Range("E1").Value = 4.734210014343262
